$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2000-2009 data rows (rows 2-11), shifting 2010-2020 up to rows 2-12.
$ws.Range("A2:D11").EntireRow.Delete()

# Carry the year-column formatting (border + bold/centered style) down onto the
# two freshly appended rows before filling in their values.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Append the new 2021 and 2022 data rows after the (now shifted) 2020 row (row 12).
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 645276
$ws.Range("C13").Value = 473170
$ws.Range("D13").Value = 341745

$ws.Range("A14").Value = "2022年"
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = 493526.01
$ws.Range("D14").ClearContents()
